# The data rows (3-25) of the sheet got reshuffled: each target row's
# record (id, taxon order, red-list code, taxon id, Swedish/scientific
# name, author, activity note, coordinates) now holds what used to be
# the contents of a different row. Every other column (validation
# status, site name, dates, observer, etc.) is identical across all of
# these rows already, so only the columns that actually vary need to
# be touched.
#
# Capture a snapshot of every source row's varying cells first (so
# overwriting row N later doesn't clobber data another target row
# still needs), then write all rows back out per the permutation below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually differ row-to-row and need to move with the record.
$cols = @("A", "B", "D", "E", "F", "G", "H", "M", "Q", "R")

# target row -> source row (both reference the *original*, pre-edit sheet)
$mapping = [ordered]@{
    3  = 4
    4  = 14
    5  = 16
    6  = 22
    7  = 24
    8  = 25
    9  = 3
    10 = 5
    11 = 6
    12 = 7
    13 = 8
    14 = 9
    15 = 10
    16 = 11
    17 = 12
    18 = 13
    19 = 15
    20 = 17
    21 = 18
    22 = 19
    23 = 20
    24 = 21
    25 = 23
}

# 1) Snapshot every involved source row's relevant cells before any writes.
$snapshots = @{}
foreach ($srcRow in $mapping.Values) {
    if (-not $snapshots.ContainsKey($srcRow)) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Range($col + $srcRow).Value2
        }
        $snapshots[$srcRow] = $rowVals
    }
}

# 2) Write each target row from its captured source-row snapshot.
foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    $rowVals = $snapshots[$srcRow]
    foreach ($col in $cols) {
        $cell = $ws.Range($col + $targetRow)
        $val = $rowVals[$col]
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value2 = $val
        }
    }
}
